$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 1.484391333333334
    "H2" = 4.453174000000001
    "I2" = 0.14812258302985
    "J2" = 0.157784627403132
    "M2" = 63.46725166666666
    "N2" = 190.401755
    "O2" = 0.2354497988808272
    "P2" = 0.2397164477183668
    "Q2" = 94.21023832448557
    "R2" = 847.89214492037
    "S2" = 0.03487543238408681
    "T2" = 0.03782357038564487
    "G3" = 1.484391333333334
    "H3" = 4.453174000000001
    "I3" = 0.14812258302985
    "J3" = 0.157784627403132
    "O3" = 0.1779985000094065
    "P3" = 0.1812240584798697
    "Q3" = 71.22232079618335
    "R3" = 641.0008871656501
    "S3" = 0.02636559759683206
    "T3" = 0.02859437054372964
    "G4" = 1.484391333333334
    "H4" = 4.453174000000001
    "I4" = 0.14812258302985
    "J4" = 0.157784627403132
    "M4" = 64.53809233333334
    "N4" = 193.614277
    "O4" = 0.2394223865221556
    "P4" = 0.243761023683841
    "Q4" = 95.79978492946647
    "R4" = 862.1980643651982
    "S4" = 0.03546386232683284
    "T4" = 0.03846174229736089
    "G5" = 1.484391333333334
    "H5" = 4.453174000000001
    "I5" = 0.14812258302985
    "J5" = 0.157784627403132
    "M5" = 14.3933435
    "N5" = 28.786687
    "O5" = 0.0533961963580272
    "P5" = 0.03624253541791403
    "Q5" = 21.36535434908967
    "R5" = 128.192126094538
    "S5" = 0.007909182528520057
    "T5" = 0.005718514947060379
    "G6" = 1.484391333333334
    "H6" = 4.453174000000001
    "I6" = 0.14812258302985
    "J6" = 0.157784627403132
    "M6" = 79.17795566666666
    "N6" = 237.533867
    "O6" = 0.2937331182295834
    "P6" = 0.2990559347000084
    "Q6" = 117.5310711826509
    "R6" = 1057.779640643858
    "S6" = 0.0435085081935782
    "T6" = 0.04718642922933618
    "I7" = 0.4232592596904365
    "J7" = 0.4508684848665014
    "M7" = 63.46725166666666
    "N7" = 190.401755
    "O7" = 0.2354497988808272
    "P7" = 0.2397164477183668
    "Q7" = 269.2051064248967
    "R7" = 2422.84595782407
    "S7" = 0.09965630756856111
    "T7" = 0.1080805915803599
    "I8" = 0.4232592596904365
    "J8" = 0.4508684848665014
    "O8" = 0.1779985000094065
    "P8" = 0.1812240584798697
    "S8" = 0.07533951333998953
    "T8" = 0.0817082166681771
    "I9" = 0.4232592596904365
    "J9" = 0.4508684848665014
    "M9" = 64.53809233333334
    "N9" = 193.614277
    "O9" = 0.2394223865221556
    "P9" = 0.243761023683841
    "Q9" = 273.7472248885754
    "R9" = 2463.725023997178
    "S9" = 0.1013377420726851
    "T9" = 0.1099041634178408
    "I10" = 0.4232592596904365
    "J10" = 0.4508684848665014
    "M10" = 14.3933435
    "N10" = 28.786687
    "O10" = 0.0533961963580272
    "P10" = 0.03624253541791403
    "Q10" = 61.051352736653
    "R10" = 366.308116419918
    "S10" = 0.02260043454078377
    "T10" = 0.01634061703159542
    "I11" = 0.4232592596904365
    "J11" = 0.4508684848665014
    "M11" = 79.17795566666666
    "N11" = 237.533867
    "O11" = 0.2937331182295834
    "P11" = 0.2990559347000084
    "Q11" = 335.8442255180486
    "R11" = 3022.598029662438
    "S11" = 0.1243252621684169
    "T11" = 0.1348348961685281
    "G12" = 0.8171586666666667
    "H12" = 2.451476
    "I12" = 0.08154160546066344
    "J12" = 0.08686056894424524
    "M12" = 63.46725166666666
    "N12" = 190.401755
    "O12" = 0.2354497988808272
    "P12" = 0.2397164477183668
    "Q12" = 51.86281474893111
    "R12" = 466.76533274038
    "S12" = 0.01919895460613297
    "T12" = 0.02082190703411076
    "G13" = 0.8171586666666667
    "H13" = 2.451476
    "I13" = 0.08154160546066344
    "J13" = 0.08686056894424524
    "O13" = 0.1779985000094065
    "P13" = 0.1812240584798697
    "Q13" = 39.20794698256667
    "R13" = 352.8715228431
    "S13" = 0.01451428346035692
    "T13" = 0.01574122482594665
    "G14" = 0.8171586666666667
    "H14" = 2.451476
    "I14" = 0.08154160546066344
    "J14" = 0.08686056894424524
    "M14" = 64.53809233333334
    "N14" = 193.614277
    "O14" = 0.2394223865221556
    "P14" = 0.243761023683841
    "Q14" = 52.73786148031689
    "R14" = 474.640753322852
    "S14" = 0.01952288578024008
    "T14" = 0.02117322120361007
    "G15" = 0.8171586666666667
    "H15" = 2.451476
    "I15" = 0.08154160546066344
    "J15" = 0.08686056894424524
    "M15" = 14.3933435
    "N15" = 28.786687
    "O15" = 0.0533961963580272
    "P15" = 0.03624253541791403
    "Q15" = 11.76164538333533
    "R15" = 70.56987230001199
    "S15" = 0.004354011576526368
    "T15" = 0.003148047246381971
    "G16" = 0.8171586666666667
    "H16" = 2.451476
    "I16" = 0.08154160546066344
    "J16" = 0.08686056894424524
    "M16" = 79.17795566666666
    "N16" = 237.533867
    "O16" = 0.2937331182295834
    "P16" = 0.2990559347000084
    "Q16" = 64.70095268196577
    "R16" = 582.3085741376919
    "S16" = 0.02395147003740709
    "T16" = 0.02597616863419578
    "G17" = 1.8409955
    "H17" = 3.681991
    "I17" = 0.1837069529326701
    "J17" = 0.1304601118296041
    "M17" = 63.46725166666666
    "N17" = 190.401755
    "O17" = 0.2354497988808272
    "P17" = 0.2397164477183668
    "Q17" = 116.8429247157008
    "R17" = 701.0575482942049
    "S17" = 0.04325376512100677
    "T17" = 0.03127343457673357
    "G18" = 1.8409955
    "H18" = 3.681991
    "I18" = 0.1837069529326701
    "J18" = 0.1304601118296041
    "O18" = 0.1779985000094065
    "P18" = 0.1812240584798697
    "Q18" = 88.3324829112875
    "R18" = 529.994897467725
    "S18" = 0.03269956206331391
    "T18" = 0.02364251093549851
    "G19" = 1.8409955
    "H19" = 3.681991
    "I19" = 0.1837069529326701
    "J19" = 0.1304601118296041
    "M19" = 64.53809233333334
    "N19" = 193.614277
    "O19" = 0.2394223865221556
    "P19" = 0.243761023683841
    "Q19" = 118.8143375642512
    "R19" = 712.8860253855071
    "S19" = 0.0439835570918532
    "T19" = 0.03180109040949267
    "G20" = 1.8409955
    "H20" = 3.681991
    "I20" = 0.1837069529326701
    "J20" = 0.1304601118296041
    "M20" = 14.3933435
    "N20" = 28.786687
    "O20" = 0.0533961963580272
    "P20" = 0.03624253541791403
    "Q20" = 26.49808061345425
    "R20" = 105.992322453817
    "S20" = 0.009809252531127714
    "T20" = 0.00472820522360945
    "G21" = 1.8409955
    "H21" = 3.681991
    "I21" = 0.1837069529326701
    "J21" = 0.1304601118296041
    "M21" = 79.17795566666666
    "N21" = 237.533867
    "O21" = 0.2937331182295834
    "P21" = 0.2990559347000084
    "Q21" = 145.7662600815328
    "R21" = 874.597560489197
    "S21" = 0.0539608161253685
    "T21" = 0.03901487068426986
    "G22" = 1.637187333333333
    "H22" = 4.911562
    "I22" = 0.1633695988863799
    "J22" = 0.1740262069565172
    "M22" = 63.46725166666666
    "N22" = 190.401755
    "O22" = 0.2354497988808272
    "P22" = 0.2397164477183668
    "Q22" = 103.9077805101455
    "R22" = 935.1700245913099
    "S22" = 0.03846533920103956
    "T22" = 0.04171694414151764
    "G23" = 1.637187333333333
    "H23" = 4.911562
    "I23" = 0.1633695988863799
    "J23" = 0.1740262069565172
    "O23" = 0.1779985000094065
    "P23" = 0.1812240584798697
    "Q23" = 78.55359893288333
    "R23" = 706.98239039595
    "S23" = 0.02907954354891402
    "T23" = 0.03153773550651778
    "G24" = 1.637187333333333
    "H24" = 4.911562
    "I24" = 0.1633695988863799
    "J24" = 0.1740262069565172
    "M24" = 64.53809233333334
    "N24" = 193.614277
    "O24" = 0.2394223865221556
    "P24" = 0.243761023683841
    "Q24" = 105.6609472856305
    "R24" = 950.9485255706741
    "S24" = 0.03911433925054437
    "T24" = 0.04242080635553661
    "G25" = 1.637187333333333
    "H25" = 4.911562
    "I25" = 0.1633695988863799
    "J25" = 0.1740262069565172
    "M25" = 14.3933435
    "N25" = 28.786687
    "O25" = 0.0533961963580272
    "P25" = 0.03624253541791403
    "Q25" = 23.56459966251567
    "R25" = 141.387597975094
    "S25" = 0.008723315181069282
    "T25" = 0.006307150969266812
    "G26" = 1.637187333333333
    "H26" = 4.911562
    "I26" = 0.1633695988863799
    "J26" = 0.1740262069565172
    "M26" = 79.17795566666666
    "N26" = 237.533867
    "O26" = 0.2937331182295834
    "P26" = 0.2990559347000084
    "Q26" = 129.6291460966949
    "R26" = 1166.662314870254
    "S26" = 0.04798706170481264
    "T26" = 0.05204356998367835
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
